$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...com grande demanda de doentes..." -> split the run and change
#    "demanda" to "afluência", producing three runs:
#    "...com grande " | "afluência" | " de doentes..."
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("demanda", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $r1.Start
$end1 = $r1.End
# Toggle Bold on the exact "demanda" span first: this forces Word to split
# the enclosing run into three runs (before / "demanda" / after) while
# keeping the surrounding text runs intact.
$r1.Bold = 1
# Re-seat a fresh Range over the same (now-isolated) span and swap its
# text for the replacement word, then clear the temporary Bold flag.
$r1b = $d.Range($start1, $end1)
$r1b.Text = "afluência"
$r1b.Bold = 0

# ---------------------------------------------------------------------
# 2) "Tal facto leva a questionar..." -> split into "Tal fa" | "c" | "to
#    leva a questionar..." (same text, now spread over three runs).
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Tal facto leva", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base2 = $r2.Start
$rc = $d.Range($base2 + 6, $base2 + 7)
$rc.Bold = 1
$rc.Bold = 0

# ---------------------------------------------------------------------
# 3) "...software de programação estatística R..." -> split into
#    "...software de programação" | " " | "estatística R..."
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("de programação estatística R. O R", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base3 = $r3.Start
$rSpace = $d.Range($base3 + 14, $base3 + 15)
$rSpace.Bold = 1
$rSpace.Bold = 0
$rEst = $d.Range($base3 + 15, $base3 + 15)
$rEst.Bold = 1
$rEst.Bold = 0

# ---------------------------------------------------------------------
# 4) Citation [2]: fix the "risck" typo to "risk" and collapse the three
#    runs (incl. the proofErr-wrapped misspelling run) back into one.
# ---------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("[2] Coyle M., Chang H., Burns P., Traynor V. Impact of Interactive Education on health care practitioners and older adults at risck of delirium: a literature review. Journal of Gerontological Nursing. 208;44(8):41-48.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4.Text = "[2] Coyle M., Chang H., Burns P., Traynor V. Impact of Interactive Education on health care practitioners and older adults at risk of delirium: a literature review. Journal of Gerontological Nursing. 208;44(8):41-48."
